# PlanetData.xlsx edit: add Eccentricity/Temperature/Velocity/Rotation columns
# for each planet, add Earth's Moon row, and remove the old scratch
# "Tests" calculation area (rows 18-23/28 helper cells), replacing it with
# a single reference link to the data source.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New header cells for columns K:N on row 1
# ---------------------------------------------------------------------
$ws.Range("K1").Value = "Eccentricity"
$ws.Range("L1").Value = "Temperature"
$ws.Range("M1").Value = "Velocity Km/s"
$ws.Range("N1").Value = "Rotation (hours)"

# ---------------------------------------------------------------------
# 2. New per-planet data in columns K (Eccentricity), M (Velocity Km/s)
#    and N (Rotation (hours)) for rows 2-10 (Sun .. Neptune)
# ---------------------------------------------------------------------
$ws.Range("K2").Value = 0

$ws.Range("K3").Value = 0.2056
$ws.Range("M3").Value = 47.4
$ws.Range("N3").Value = 1407.6

$ws.Range("K4").Value = 0.0068
$ws.Range("M4").Value = 35
$ws.Range("N4").Value = -5835.5

$ws.Range("K5").Value = 0.0167
$ws.Range("M5").Value = 29.8
$ws.Range("N5").Value = 23.9

$ws.Range("K6").Value = 0.0934
$ws.Range("M6").Value = 24.1
$ws.Range("N6").Value = 24.6

$ws.Range("K7").Value = 0.0484
$ws.Range("M7").Value = 13.1
$ws.Range("N7").Value = 9.9

$ws.Range("K8").Value = 0.0542
$ws.Range("M8").Value = 9.7
$ws.Range("N8").Value = 10.9

$ws.Range("K9").Value = 0.0472
$ws.Range("M9").Value = 6.8
$ws.Range("N9").Value = -17.2

$ws.Range("K10").Value = 0.0086
$ws.Range("M10").Value = 5.4
$ws.Range("N10").Value = 16.1

# ---------------------------------------------------------------------
# 3. New row 12: Earth's Moon
# ---------------------------------------------------------------------
$ws.Range("B12").Value = "Earths Moon"
$ws.Range("C12").Value = 0.073
$ws.Range("D12").Value = 3.34
$ws.Range("E12").Value = 1737.5
$ws.Range("F12").Value = 0.378
$ws.Range("G12").Value = 0.3633
$ws.Range("H12").Value = 0.3844
$ws.Range("I12").Value = 0.3844
$ws.Range("K12").Value = 0.0549

# ---------------------------------------------------------------------
# 4. Remove the old "Tests" scratch-work block (rows 18-23) entirely
# ---------------------------------------------------------------------
$ws.Range("A18").ClearContents()

$ws.Range("A19").ClearContents()
$ws.Range("C19").ClearContents()
$ws.Range("D19").ClearContents()
$ws.Range("E19").ClearContents()
$ws.Range("F19").ClearContents()

$ws.Range("A20").ClearContents()

$ws.Range("C21").ClearContents()
$ws.Range("D21").ClearContents()
$ws.Range("E21").ClearContents()

$ws.Range("A22").ClearContents()
$ws.Range("C22").ClearContents()
$ws.Range("D22").ClearContents()
$ws.Range("E22").ClearContents()
$ws.Range("F22").ClearContents()

$ws.Range("A23").ClearContents()

$ws.Range("A24").ClearContents()
$ws.Range("B24").ClearContents()

$ws.Range("E25").ClearContents()

$ws.Range("F28").ClearContents()
$ws.Range("G28").ClearContents()
$ws.Range("H28").ClearContents()

# ---------------------------------------------------------------------
# 5. Row 20 now just holds the source-data link
# ---------------------------------------------------------------------
$ws.Range("B20").Value = "http://nssdc.gsfc.nasa.gov/planetary/factsheet/moonfact.html"

# ---------------------------------------------------------------------
# 6. Column widths for the new columns K:N
# ---------------------------------------------------------------------
$ws.Columns("K").ColumnWidth = 18.26
$ws.Columns("L").ColumnWidth = 15.09
$ws.Columns("M").ColumnWidth = 13.97
$ws.Columns("N").ColumnWidth = 15.26

# ---------------------------------------------------------------------
# 7. Selection / scroll position
# ---------------------------------------------------------------------
$ws.Range("K12").Select()
